# Update "height and weight jbrelsf2" — refresh the CIBMTR ValueSet metadata
# sheet (Version/Status/Date/Contact/Jurisdiction) to match the new IG build.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# --- Simple value updates -------------------------------------------------
$ws1.Range("B3").Value2  = "0.1.7"
$ws1.Range("B6").Value2  = "draft"
$ws1.Range("B8").Value2  = "2024-11-22T12:33:30-06:00"
$ws1.Range("B10").Value2 = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Make room for two new rows: a second "Contact" row (the person) and a
# "Jurisdiction" row. Row 11 already duplicates row 10 (Contact) and is being
# replaced outright; rows 12-15 (Description/Purpose/Copyright/Immutable)
# shift down by one, to 13-16, by copying values bottom-up.
for ($r = 15; $r -ge 12; $r--) {
    $src = $r
    $dst = $r + 1
    $ws1.Range("A" + $dst).Value2 = $ws1.Range("A" + $src).Value2
    $ws1.Range("B" + $dst).Value2 = $ws1.Range("B" + $src).Value2
}

# New row 11: second Contact (person)
$ws1.Range("A11").Value2 = "Contact"
$ws1.Range("B11").Value2 = "Bob Milius (bmilius@nmdp.org)"

# New row 12: Jurisdiction (no value supplied)
$ws1.Range("A12").Value2 = "Jurisdiction"
$ws1.Range("B12").Value2 = ""

# Rows 13-16 now hold what used to be in 12-15 (Description/Purpose/Copyright/
# Immutable); row 16 is brand new territory (sheet used to stop at row 15), so
# copy the established "data row" formatting (border + top/wrap alignment)
# down onto it to match the rest of the table.
$ws1.Range("A15:B15").Copy()
$ws1.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "Metadata table updated"
